# Revert to original OR files where customized to state
# Updates the "About" sheet source citation from the California ARB LCFS
# reference to the Oregon DEQ Clean Fuels Program reference, and removes
# the now-redundant "(Boolean)" header label on the "BVTStL" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Source block (rows 3-7)
$about.Range("B3").Value = "Oregon Department of Environmental Quality"
$about.Range("B4").Value = "undated"
$about.Range("B5").Value = "Oregon Clean Fuels Program: Exemptions"
$about.Range("B6").Value = "https://secure.sos.state.or.us/oard/view.action?ruleNumber=340-253-0250"
$about.Range("B7").Value = "Section (2)"

# Highlight/box the primary source-organization cell to match the new
# "source" callout styling.
$srcCell = $about.Range("B3")
$srcCell.Font.Bold = $true
$srcCell.Interior.Color = 12566463
$srcCell.Borders.Weight = -4138

# Notes block (rows 9-13) condensed to three lines (rows 10-12)
$about.Range("A10").Value = "The Oregon CFP identifies aircrafts as exempt. It also"
$about.Range("A11").Value = "identifies watercraft, however we assume these are smaller"
$about.Range("A12").Value = "than ships here."
$about.Range("A13").ClearContents()

# Old standalone note (row 15) removed entirely
$about.Range("A15").ClearContents()

# ---------------------------------------------------------------------
# "BVTStL" sheet
# ---------------------------------------------------------------------
$bvtstl = $wb.Worksheets.Item("BVTStL")

# Drop the now-unused "(Boolean)" header label above the row labels
$bvtstl.Range("A1").ClearContents()
